$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "datos actualizados" timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 28 de Mayo de 2020 a las 12:10"

# --- Country reordering + refreshed statistics ---
# Row 13
$ws.Range("B13").Value = 158897
$ws.Range("C13").Value = 811
$ws.Range("D13").Value = 67901
$ws.Range("E13").Value = 86456

# Row 62
$ws.Range("A62").Value = "Marruecos"
$ws.Range("B62").Value = 7636
$ws.Range("C62").Value = 35
$ws.Range("D62").Value = 5109
$ws.Range("E62").Value = 2325
$ws.Range("H62").Value = 202

# Row 63
$ws.Range("A63").Value = "Malasia"
$ws.Range("B63").Value = 7629
$ws.Range("C63").Value = 10
$ws.Range("D63").Value = 6169
$ws.Range("E63").Value = 1345
$ws.Range("H63").Value = 115

# Row 110
$ws.Range("A110").Value = "Tunez"
$ws.Range("B110").Value = 1068
$ws.Range("C110").Value = 17
$ws.Range("D110").Value = 938
$ws.Range("E110").Value = 82
$ws.Range("H110").Value = 48

# Row 111
$ws.Range("A111").Value = "Hong Kong"
$ws.Range("B111").Value = 1067
$ws.Range("C111").Value = 0
$ws.Range("D111").Value = 1035
$ws.Range("E111").Value = 28
$ws.Range("G111").Value = 0
$ws.Range("H111").Value = 4

# Row 112
$ws.Range("A112").Value = "Letonia"
$ws.Range("B112").Value = 1061
$ws.Range("C112").Value = 4
$ws.Range("D112").Value = 741
$ws.Range("E112").Value = 296
$ws.Range("G112").Value = 1
$ws.Range("H112").Value = 24

# Row 113
$ws.Range("A113").Value = "Zambia"
$ws.Range("B113").Value = 1057
$ws.Range("D113").Value = 779
$ws.Range("E113").Value = 271
$ws.Range("H113").Value = 7

# Row 200
$ws.Range("A200").Value = "Nueva Caledonia"
$ws.Range("D200").Value = 18
$ws.Range("H200").Value = 0

# Row 201
$ws.Range("A201").Value = "Belice"
$ws.Range("D201").Value = 16
$ws.Range("H201").Value = 2

# Row 213
$ws.Range("A213").Value = "Islas Virgenes Britanicas"
$ws.Range("D213").Value = 7
$ws.Range("H213").Value = 1

# Row 214
$ws.Range("A214").Value = "Papua Nueva Guinea"
$ws.Range("D214").Value = 8
$ws.Range("H214").Value = 0

# Row 215
$ws.Range("A215").Value = "Bonaire, San Eustaquio y Saba"

# Row 216
$ws.Range("A216").Value = "San Bartolome"
